# APITesting generated file - 2023-11-24 12:03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Request Header changes from None to ContentType.JSON
$ws.Range("E4").Value = "ContentType.JSON"

# Row 6: Assertions changes from None to a data-array assertion
$ws.Range("N6").Value = "data array should not be null, should have more than 1 item"

# Row 7: endpoint/auth test now targets /register with BasicAuth (expects 400)
$ws.Range("C7").Value = "/register"
$ws.Range("D7").Value = "POST"
$ws.Range("H7").Value = "None"
$ws.Range("J7").Value = "BasicAuth"
$ws.Range("K7").Value = "{username: 'testuser', password: 'testpass'}"
$ws.Range("L7").Value = "400"
$ws.Range("N7").Value = "response should contain error message"

# Row 8: now a GET of a single user expected to 404
$ws.Range("C8").Value = "/users/5"
$ws.Range("D8").Value = "GET"
$ws.Range("H8").Value = "id=5"
$ws.Range("J8").Value = "None"
$ws.Range("K8").Value = "None"
$ws.Range("L8").Value = "404"
$ws.Range("N8").Value = "None"

# Row 9: new row - GET /users?page=3
$ws.Range("A9").Value = "Rest Assured"
$ws.Range("B9").Value = "https://reqres.in/api"
$ws.Range("C9").Value = "/users?page=3"
$ws.Range("D9").Value = "GET"
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "None"
$ws.Range("G9").Value = "None"
$ws.Range("H9").Value = "None"
$ws.Range("I9").Value = "None"
$ws.Range("J9").Value = "None"
$ws.Range("K9").Value = "None"
$ws.Range("L9").Value = "200"
$ws.Range("M9").Value = "None"
$ws.Range("N9").Value = "data array should not be null, should have more than 1 item. Page number should be 3."
